$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Tracks" test-data rows (57, 59, 61) -----------------------------
# Row 57: mirrors the other "selection tab" test rows (s / test)
$ws.Range("B57").Value = "s"
$ws.Range("C57").Value = "test"

# Row 59: "Add Track" test data
$ws.Range("B59").Value = "test track"
$ws.Range("C59").Value = "test"
$ws.Range("D59").Value = "this is  a  test track for aut."

# Row 61: "Edit Track" test data
$ws.Range("B61").Value = "test"
$ws.Range("C61").Value = "track test"
$ws.Range("D61").Value = "this is  a  test track for aut."
$ws.Range("E61").Value = "ttc"
$ws.Range("F61").Value = "trial test course"

# --- Relabel the "Quick sales tips" add/edit headers ----------------------
$ws.Range("A34").Value = "add sales tips"
$ws.Range("A36").Value = "edit sales tips"

# --- Update the view / selection state ------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("L31").Select()
